$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "_duplicates" column (M) entirely - the whole column is removed.
$ws.Columns.Item(13).Delete()

# Rename header labels to match the new Dataset-object naming scheme
# (breaking rename: "_link"/"_merge"/etc suffixes become "_x"/"_mp_*").
$ws.Range("J1").Value = "_mp_merge"
$ws.Range("K1").Value = "_mp_diff_days"
$ws.Range("L1").Value = "_mp_abs_diff_days"
$ws.Range("A1").Value = "PIDN_x"
$ws.Range("B1").Value = "DCDate_x"

# Re-fit the columns whose header text grew so they keep showing fully.
$ws.Columns.Item(10).AutoFit()
$ws.Columns.Item(11).AutoFit()
$ws.Columns.Item(12).AutoFit()
